$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct tiny precision drift on the last existing row's timestamp
$ws.Range("A16").Value = 45816.39137637732

# Append the new price row with the same date/number formatting as the
# other rows in column A
$ws.Range("A17").Value = 45817.39392137039
$ws.Range("A17").NumberFormat = $ws.Range("A16").NumberFormat

$ws.Range("B17").Value = "EVOWHEY PROTEIN"
$ws.Range("C17").Value = "2Kg"
$ws.Range("D17").Value = "37,90€"
